# "ADD threads for CRUD BBDD" - refresh the exported stock data:
# insert a new "Megaman 2" / NES row, append two more rows at the
# bottom, and replace the broken Java object-reference placeholders
# ("[B@...]" / blank) shown in the Imagen column with readable text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new "Megaman 2" row right before the old row 15 (Metal Slug 3)
$ws.Rows.Item(15).Insert()

# Insert 2 more rows just above the last row (Megaman), which has shifted to row 26
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).Insert()

# Rewrite every data row (2-28) so values/shifts and the Imagen placeholder text are correct
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "World of Warcraft"
$ws.Cells.Item(2,3).Value = "Blizzard"
$ws.Cells.Item(2,4).Value = 45
$ws.Cells.Item(2,5).Value = 90
$ws.Cells.Item(2,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(3,1).Value = 9
$ws.Cells.Item(3,2).Value = "Ghost & Goblins"
$ws.Cells.Item(3,3).Value = "NINTENDO"
$ws.Cells.Item(3,4).Value = 89
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(4,1).Value = 86
$ws.Cells.Item(4,2).Value = "Sonic & Knuckles"
$ws.Cells.Item(4,3).Value = "SEGA"
$ws.Cells.Item(4,4).Value = 87
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(5,1).Value = 154
$ws.Cells.Item(5,2).Value = "ZELDA OCARINA OF TIME"
$ws.Cells.Item(5,3).Value = "SWITCH"
$ws.Cells.Item(5,4).Value = 55
$ws.Cells.Item(5,5).Value = 4
$ws.Cells.Item(5,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(6,1).Value = 157
$ws.Cells.Item(6,2).Value = "Mario Bros 3"
$ws.Cells.Item(6,3).Value = "Nintendo NES"
$ws.Cells.Item(6,4).Value = 1500
$ws.Cells.Item(6,5).Value = 10
$ws.Cells.Item(6,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(7,1).Value = 347
$ws.Cells.Item(7,2).Value = "Zelda Remastered"
$ws.Cells.Item(7,3).Value = "Nintendo Switch"
$ws.Cells.Item(7,4).Value = 500
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(8,1).Value = 547
$ws.Cells.Item(8,2).Value = "Super Mario Bros 3"
$ws.Cells.Item(8,3).Value = "Nintendo NES"
$ws.Cells.Item(8,4).Value = 1000
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(9,1).Value = 554
$ws.Cells.Item(9,2).Value = "Tomb Raider"
$ws.Cells.Item(9,3).Value = "SONY"
$ws.Cells.Item(9,4).Value = 87
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(10,1).Value = 555
$ws.Cells.Item(10,2).Value = "Gears of war"
$ws.Cells.Item(10,3).Value = "XBOX"
$ws.Cells.Item(10,4).Value = 65
$ws.Cells.Item(10,5).Value = 4
$ws.Cells.Item(10,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(11,1).Value = 556
$ws.Cells.Item(11,2).Value = "Gears of war 2"
$ws.Cells.Item(11,3).Value = "XBOX"
$ws.Cells.Item(11,4).Value = 65
$ws.Cells.Item(11,5).Value = 4
$ws.Cells.Item(11,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(12,1).Value = 4467
$ws.Cells.Item(12,2).Value = "CS2 Global Offensive"
$ws.Cells.Item(12,3).Value = "Steam"
$ws.Cells.Item(12,4).Value = 15.5
$ws.Cells.Item(12,5).Value = 5
$ws.Cells.Item(12,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(13,1).Value = 9002
$ws.Cells.Item(13,2).Value = "Punch Out"
$ws.Cells.Item(13,3).Value = "NINTENDO"
$ws.Cells.Item(13,4).Value = 54
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(14,1).Value = 9003
$ws.Cells.Item(14,2).Value = "Punch Out 2"
$ws.Cells.Item(14,3).Value = "NINTENDO"
$ws.Cells.Item(14,4).Value = 54
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(15,1).Value = 9986
$ws.Cells.Item(15,2).Value = "Megaman 2"
$ws.Cells.Item(15,3).Value = "NES"
$ws.Cells.Item(15,4).Value = 160
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = "javax.swing.ImageIcon@115b973f"

$ws.Cells.Item(16,1).Value = 10101
$ws.Cells.Item(16,2).Value = "Metal Slug 3"
$ws.Cells.Item(16,3).Value = "SNK"
$ws.Cells.Item(16,4).Value = 99999
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(17,1).Value = 24024
$ws.Cells.Item(17,2).Value = "Pac-Man"
$ws.Cells.Item(17,3).Value = "NES"
$ws.Cells.Item(17,4).Value = 50
$ws.Cells.Item(17,5).Value = 4
$ws.Cells.Item(17,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(18,1).Value = 45484
$ws.Cells.Item(18,2).Value = "ZELDA A LINK TO THE PAST"
$ws.Cells.Item(18,3).Value = "NINTENDO 64"
$ws.Cells.Item(18,4).Value = 150
$ws.Cells.Item(18,5).Value = 5
$ws.Cells.Item(18,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(19,1).Value = 88965
$ws.Cells.Item(19,2).Value = "Zelda A Link To The Past"
$ws.Cells.Item(19,3).Value = "NES"
$ws.Cells.Item(19,4).Value = 97
$ws.Cells.Item(19,5).Value = 44
$ws.Cells.Item(19,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(20,1).Value = 90909
$ws.Cells.Item(20,2).Value = "Super Mario Brosh Land 3"
$ws.Cells.Item(20,3).Value = "NINTENDO"
$ws.Cells.Item(20,4).Value = 157
$ws.Cells.Item(20,5).Value = 4
$ws.Cells.Item(20,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(21,1).Value = 99984
$ws.Cells.Item(21,2).Value = "Mario Bros 1"
$ws.Cells.Item(21,3).Value = "NINTENDO"
$ws.Cells.Item(21,4).Value = 89
$ws.Cells.Item(21,5).Value = 1
$ws.Cells.Item(21,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(22,1).Value = 121212
$ws.Cells.Item(22,2).Value = "Mario Bros Land"
$ws.Cells.Item(22,3).Value = "NES"
$ws.Cells.Item(22,4).Value = 8888
$ws.Cells.Item(22,5).Value = 9
$ws.Cells.Item(22,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(23,1).Value = 987897
$ws.Cells.Item(23,2).Value = "Zelda Ocarina"
$ws.Cells.Item(23,3).Value = "NES"
$ws.Cells.Item(23,4).Value = 88
$ws.Cells.Item(23,5).Value = 3
$ws.Cells.Item(23,6).Value = "javax.swing.ImageIcon@6a6e279b"

$ws.Cells.Item(24,1).Value = 9875555
$ws.Cells.Item(24,2).Value = "Super Mario Bros 3. Ed Limitada"
$ws.Cells.Item(24,3).Value = "NINTENDO"
$ws.Cells.Item(24,4).Value = 1600
$ws.Cells.Item(24,5).Value = 1
$ws.Cells.Item(24,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(25,1).Value = 48592378
$ws.Cells.Item(25,2).Value = "Mario Bros"
$ws.Cells.Item(25,3).Value = "NES"
$ws.Cells.Item(25,4).Value = 45
$ws.Cells.Item(25,5).Value = 4
$ws.Cells.Item(25,6).Value = "javax.swing.ImageIcon@17ace130"

$ws.Cells.Item(26,1).Value = 379827496
$ws.Cells.Item(26,2).Value = "Megaman 2"
$ws.Cells.Item(26,3).Value = "Nes"
$ws.Cells.Item(26,4).Value = 88
$ws.Cells.Item(26,5).Value = 3
$ws.Cells.Item(26,6).Value = "javax.swing.ImageIcon@3b203b2a"

$ws.Cells.Item(27,1).Value = 379827498
$ws.Cells.Item(27,2).Value = "Megaman 2"
$ws.Cells.Item(27,3).Value = "Nes"
$ws.Cells.Item(27,4).Value = 88
$ws.Cells.Item(27,5).Value = 3
$ws.Cells.Item(27,6).Value = "IMAGEN NO DISPONIBLE"

$ws.Cells.Item(28,1).Value = 797686286
$ws.Cells.Item(28,2).Value = "Megaman"
$ws.Cells.Item(28,3).Value = "NES"
$ws.Cells.Item(28,4).Value = 150
$ws.Cells.Item(28,5).Value = 3
$ws.Cells.Item(28,6).Value = "javax.swing.ImageIcon@5463bdca"

